# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 4-6 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1724
$wsExhibit.Range("F5").Value = 768
$wsExhibit.Range("F6").Value = 192

# Sheet "全部类型": rows 4, 6, 7 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1724
$wsAll.Range("F6").Value = 768
$wsAll.Range("F7").Value = 192
